$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30
$ws.Range("A30").Value = "POINT (6.9355980807587105 49.364007303777186)"
$ws.Range("B30").Value = 10043
$ws.Range("C30").Value = 364882131
$ws.Range("D30").Value = 300
$ws.Range("E30").Value = 49.36400730377719
$ws.Range("F30").Value = 6.935598080758711

# Row 31
$ws.Range("A31").Value = "POINT (7.001045558866386 49.401271550367575)"
$ws.Range("B31").Value = 10043
$ws.Range("C31").Value = 419910303
$ws.Range("D31").Value = 300
$ws.Range("E31").Value = 49.40127155036758
$ws.Range("F31").Value = 7.001045558866386

# Row 32
$ws.Range("A32").Value = "POINT (7.029496440077044 49.403660960617195)"
$ws.Range("B32").Value = 10043
$ws.Range("C32").Value = 388230747
$ws.Range("D32").Value = 200
$ws.Range("E32").Value = 49.4036609606172
$ws.Range("F32").Value = 7.029496440077044

# Row 33
$ws.Range("A33").Value = "POINT (7.091227930374576 49.38422368213677)"
$ws.Range("B33").Value = 10043
$ws.Range("C33").Value = 1264098455
$ws.Range("D33").Value = 300
$ws.Range("E33").Value = 49.38422368213677
$ws.Range("F33").Value = 7.091227930374576
